$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HOT_WATER")

# Insert a new column before the existing "Tsww0_C" column (column C),
# pushing "Tsww0_C"/"Qwwmax_Wm2" one slot to the right.
$ws.Columns.Item(3).Insert()

# Give the new column the same look as its "code" neighbour, and restyle
# the shifted "Tsww0_C" column to match as well (removing the old
# dedicated numeric styling/restriction on that column).
$ws.Range("B2:B6").Copy()
$ws.Range("C2:D6").PasteSpecial(-4122)

# Match the new column's width to the neighbouring "code" column.
$ws.Columns.Item(3).ColumnWidth = $ws.Columns.Item(2).ColumnWidth()

# Header for the new "class_dhw" column.
$ws.Cells.Item(1, 3).Value = "class_dhw"

# Classification values per DHW assembly row.
$ws.Cells.Item(2, 3).Value = "NONE"
$ws.Cells.Item(3, 3).Value = "HIGH_TEMP"
$ws.Cells.Item(4, 3).Value = "MEDIUM_TEMP"
$ws.Cells.Item(5, 3).Value = "LOW_TEMP"
$ws.Cells.Item(6, 3).Value = "HIGH_TEMP"
